$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.084.01'
$ws.Range("E2").Value = '  +5.54%  '

# Row 3
$ws.Range("D3").Value = '1.921.14'
$ws.Range("E3").Value = '  +2.45%  '

# Row 4
$ws.Range("E4").Value = '  -0.86%  '

# Row 5
$ws.Range("D5").Value = "'330.88"
$ws.Range("E5").Value = '  +4.74%  '

# Row 6
$ws.Range("E6").Value = '  -0.81%  '

# Row 7
$ws.Range("D7").Value = "'0.5239"
$ws.Range("E7").Value = '  +2.96%  '

# Row 8
$ws.Range("D8").Value = "'0.4050"
$ws.Range("E8").Value = '  +3.76%  '

# Row 9
$ws.Range("D9").Value = "'0.08490"
$ws.Range("E9").Value = '  +1.48%  '

# Row 10
$ws.Range("D10").Value = "'42.96"
$ws.Range("E10").Value = '  +3.31%  '

# Row 11
$ws.Range("D11").Value = "'1.128"
$ws.Range("E11").Value = '  +2.23%  '

# Row 12
$ws.Range("D12").Value = "'22.36"
$ws.Range("E12").Value = '  +9.74%  '

# Row 13
$ws.Range("D13").Value = "'6.400"
$ws.Range("E13").Value = '  +2.79%  '

# Row 14
$ws.Range("D14").Value = '1.919.77'
$ws.Range("E14").Value = '  +2.39%  '

# Row 15
$ws.Range("D15").Value = "'7.401"
$ws.Range("E15").Value = '  +1.86%  '

# Row 16
$ws.Range("E16").Value = '  -0.94%  '

# Row 17
$ws.Range("E17").Value = '  +5.55%  '

# Row 18
$ws.Range("D18").Value = "'0.00001115"
$ws.Range("E18").Value = '  +1.01%  '

# Row 19
$ws.Range("D19").Value = "'0.06705"
$ws.Range("E19").Value = '  -0.35%  '

# Row 20
$ws.Range("D20").Value = "'18.27"
$ws.Range("E20").Value = '  +3.24%  '

# Row 21
$ws.Range("E21").Value = '  -0.74%  '

# Row 22
$ws.Range("D22").Value = "'6.069"
$ws.Range("E22").Value = '  +2.52%  '

# Row 23
$ws.Range("D23").Value = '30.089.71'
$ws.Range("E23").Value = '  +5.50%  '

# Row 24
$ws.Range("E24").Value = '  +1.37%  '

# Row 25
$ws.Range("D25").Value = "'2.223"
$ws.Range("E25").Value = '  -0.27%  '

# Row 26
$ws.Range("D26").Value = '2.141.06'
$ws.Range("E26").Value = '  +2.54%  '

# Row 27
$ws.Range("D27").Value = "'21.15"
$ws.Range("E27").Value = '  +2.64%  '

# Row 28
$ws.Range("D28").Value = "'160.08"
$ws.Range("E28").Value = '  -1.13%  '

# Row 29
$ws.Range("D29").Value = "'2.452"
$ws.Range("E29").Value = '  +2.77%  '

# Row 30
$ws.Range("D30").Value = "'129.51"
$ws.Range("E30").Value = '  +3.08%  '

# Row 31
$ws.Range("E31").Value = '  +4.07%  '

# Row 32
$ws.Range("E32").Value = '  +1.54%  '

# Row 33
$ws.Range("D33").Value = "'6.115"
$ws.Range("E33").Value = '  +6.03%  '

# Row 34
$ws.Range("D34").Value = "'3.642"
$ws.Range("E34").Value = '  +0.82%  '

# Row 35
$ws.Range("D35").Value = "'0.02522"
$ws.Range("E35").Value = '  +2.64%  '

# Row 36
$ws.Range("D36").Value = "'0.06615"
$ws.Range("E36").Value = '  +1.32%  '

# Row 37
$ws.Range("D37").Value = "'0.2230"
$ws.Range("E37").Value = '  +3.22%  '

# Row 38
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = "'1.239"
$ws.Range("E38").Value = '  +4.09%  '

# Row 39
$ws.Range("D39").Value = "'9.044"
$ws.Range("E39").Value = '  +2.56%  '

# Row 40
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = "'5.219"
$ws.Range("E40").Value = '  +3.35%  '

# Row 41
$ws.Range("D41").Value = "'0.6565"
$ws.Range("E41").Value = '  +2.72%  '

# Row 42
$ws.Range("D42").Value = "'11.69"
$ws.Range("E42").Value = '  +5.43%  '

# Row 43
$ws.Range("D43").Value = "'1.244"
$ws.Range("E43").Value = '  +0.37%  '

# Row 44
$ws.Range("E44").Value = '  +3.36%  '

# Row 45
$ws.Range("D45").Value = "'13.22"
$ws.Range("E45").Value = '  +1.55%  '

# Row 46
$ws.Range("D46").Value = "'3.786"
$ws.Range("E46").Value = '  +2.78%  '

# Row 47
$ws.Range("D47").Value = "'2.091"
$ws.Range("E47").Value = '  +4.23%  '

# Row 48
$ws.Range("D48").Value = "'1.245"
$ws.Range("E48").Value = '  +2.50%  '

# Row 49
$ws.Range("E49").Value = '  +3.29%  '

# Row 50
$ws.Range("D50").Value = "'80.24"
$ws.Range("E50").Value = '  +5.08%  '

# Row 51
$ws.Range("E51").Value = '  +0.85%  '
